# Update weekly Fruta/Hortaliza price data (Femacal de La Calera - Caqui)
# Columns updated per row: D (Fecha), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# R (Origen), S (Precio $/Kg)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; Fecha = 44333; Calidad = "Especial"; Volumen = 58; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 3; Fecha = 44333; Calidad = "Primera"; Volumen = 65; PrecioMin = 9000; PrecioMax = 9000; PrecioProm = 9000; Origen = "Provincia de Quillota"; PrecioKg = 900 }
    @{ Row = 4; Fecha = 44333; Calidad = "Segunda"; Volumen = 60; PrecioMin = 8000; PrecioMax = 8000; PrecioProm = 8000; Origen = "Provincia de Quillota"; PrecioKg = 800 }
    @{ Row = 5; Fecha = 44309; Calidad = "Primera"; Volumen = 45; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 6; Fecha = 44321; Calidad = "Primera"; Volumen = 58; PrecioMin = 9000; PrecioMax = 9000; PrecioProm = 9000; Origen = "Provincia de Quillota"; PrecioKg = 900 }
    @{ Row = 7; Fecha = 44306; Calidad = "Primera"; Volumen = 45; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 8; Fecha = 44307; Calidad = "Primera"; Volumen = 40; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 9; Fecha = 44308; Calidad = "Primera"; Volumen = 45; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 10; Fecha = 44308; Calidad = "Segunda"; Volumen = 48; PrecioMin = 8000; PrecioMax = 8000; PrecioProm = 8000; Origen = "Provincia de Quillota"; PrecioKg = 800 }
    @{ Row = 11; Fecha = 44301; Calidad = "Primera"; Volumen = 45; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 12; Fecha = 44314; Calidad = "Primera"; Volumen = 47; PrecioMin = 9000; PrecioMax = 9000; PrecioProm = 9000; Origen = "Provincia de Quillota"; PrecioKg = 900 }
    @{ Row = 13; Fecha = 44328; Calidad = "Primera"; Volumen = 45; PrecioMin = 8000; PrecioMax = 8000; PrecioProm = 8000; Origen = "Provincia de Quillota"; PrecioKg = 800 }
    @{ Row = 14; Fecha = 44328; Calidad = "Segunda"; Volumen = 48; PrecioMin = 7000; PrecioMax = 7000; PrecioProm = 7000; Origen = "Provincia de Quillota"; PrecioKg = 700 }
    @{ Row = 15; Fecha = 44319; Calidad = "Primera"; Volumen = 68; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 16; Fecha = 44319; Calidad = "Segunda"; Volumen = 57; PrecioMin = 8000; PrecioMax = 8000; PrecioProm = 8000; Origen = "Provincia de Quillota"; PrecioKg = 800 }
    @{ Row = 17; Fecha = 44329; Calidad = "Primera"; Volumen = 56; PrecioMin = 9000; PrecioMax = 9000; PrecioProm = 9000; Origen = "Región Metropolitana"; PrecioKg = 900 }
    @{ Row = 18; Fecha = 44329; Calidad = "Segunda"; Volumen = 50; PrecioMin = 8000; PrecioMax = 8000; PrecioProm = 8000; Origen = "Región Metropolitana"; PrecioKg = 800 }
    @{ Row = 19; Fecha = 44302; Calidad = "Primera"; Volumen = 45; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 20; Fecha = 44312; Calidad = "Primera"; Volumen = 48; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 21; Fecha = 44322; Calidad = "Primera"; Volumen = 56; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 22; Fecha = 44322; Calidad = "Segunda"; Volumen = 40; PrecioMin = 8000; PrecioMax = 8000; PrecioProm = 8000; Origen = "Provincia de Quillota"; PrecioKg = 800 }
    @{ Row = 23; Fecha = 44699; Calidad = "Especial"; Volumen = 56; PrecioMin = 12000; PrecioMax = 12000; PrecioProm = 12000; Origen = "Provincia de Quillota"; PrecioKg = 1200 }
    @{ Row = 24; Fecha = 44699; Calidad = "Primera"; Volumen = 60; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 25; Fecha = 44323; Calidad = "Primera"; Volumen = 60; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 26; Fecha = 44323; Calidad = "Segunda"; Volumen = 50; PrecioMin = 9000; PrecioMax = 9000; PrecioProm = 9000; Origen = "Provincia de Quillota"; PrecioKg = 900 }
    @{ Row = 27; Fecha = 44315; Calidad = "Primera"; Volumen = 45; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 28; Fecha = 44326; Calidad = "Primera"; Volumen = 65; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Provincia de Quillota"; PrecioKg = 1000 }
    @{ Row = 29; Fecha = 44326; Calidad = "Segunda"; Volumen = 67; PrecioMin = 8000; PrecioMax = 8000; PrecioProm = 8000; Origen = "Provincia de Quillota"; PrecioKg = 800 }
    @{ Row = 30; Fecha = 44343; Calidad = "Especial"; Volumen = 47; PrecioMin = 10000; PrecioMax = 10000; PrecioProm = 10000; Origen = "Región Metropolitana"; PrecioKg = 1000 }
    @{ Row = 31; Fecha = 44343; Calidad = "Primera"; Volumen = 50; PrecioMin = 9000; PrecioMax = 9000; PrecioProm = 9000; Origen = "Región Metropolitana"; PrecioKg = 900 }
    @{ Row = 32; Fecha = 44343; Calidad = "Segunda"; Volumen = 58; PrecioMin = 8000; PrecioMax = 8000; PrecioProm = 8000; Origen = "Región Metropolitana"; PrecioKg = 800 }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Cells.Item($r, 4).Value = $row.Fecha        # D: Fecha
    $ws.Cells.Item($r, 12).Value = $row.Calidad     # L: Calidad
    $ws.Cells.Item($r, 13).Value = $row.Volumen     # M: Volumen
    $ws.Cells.Item($r, 14).Value = $row.PrecioMin   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $row.PrecioMax   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $row.PrecioProm  # P: Precio promedio ponderado
    $ws.Cells.Item($r, 18).Value = $row.Origen      # R: Origen
    $ws.Cells.Item($r, 19).Value = $row.PrecioKg    # S: Precio $/Kg
}

